$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Target cell: B11 must hold the literal TEXT "1" (a shared string), not the
# number 1, while keeping its existing cell style (s="23") untouched.
#
# A plain `.Value = "1"` assignment gets auto-coerced to a number by the
# input parser. Pre-formatting B11 itself as Text (NumberFormat = "@")
# does store it as text, but it permanently mutates/replaces B11's style
# slot. So instead: stage the text in a far-away scratch cell formatted as
# Text, Copy it, and PasteSpecial *values only* into B11 - this carries
# over the "text" cell type/content without touching B11's style index.
$scratch = $ws.Range("ZZ1000")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$target = $ws.Range("B11")
$target.PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
$scratch.Clear()

